# Update the "想去人数" (interested-count) column F values on the sheets
# that list the 漫展 events: "展览" and "全部类型". Both sheets carry the
# same 23 event rows (rows 2-23) with identical F-column figures.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 380
    3  = 10884
    5  = 987
    6  = 201
    7  = 1349
    8  = 8345
    12 = 226
    14 = 3336
    18 = 837
    20 = 1079
    22 = 133
    23 = 1861
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
